$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the two abstract class headers (UML stereotype notation) ---
$ws.Range("E8").Value = "<<abstract>> Scuola "
$ws.Range("E23").Value = "<<abstract>> ScuolaSuperiore"

# --- Scuola (abstract) class: calcoloContributo() now returns double instead of int ---
$ws.Range("E18").Value = "'+ calcoloContributo(): double"

# --- ScuolaSuperiore (abstract) class gains three private attributes ---
$ws.Range("E25").Value = "'- contributoClasse: double"
$ws.Range("E26").Value = "'- contributoLaboratorio: double"
$ws.Range("E27").Value = "'- contributoSedeAggiuntiva: double"

# --- ScuolaElementare gains calcoloContributo()/toString() methods ---
$ws.Range("J10").Value = "'+ calcoloContributo(): double"
$ws.Range("J11").Value = "'+ toString(): String"

# --- ScuolaMedia gains calcoloContributo()/toString() methods ---
$ws.Range("J17").Value = "'+ calcoloContributo(): double"
$ws.Range("J18").Value = "'+ toString(): String"

# --- Liceo gains calcoloContributo()/toString() methods ---
$ws.Range("J25").Value = "'+ calcoloContributo(): double"
$ws.Range("J26").Value = "'+ toString(): String"

# --- Tecnico gains calcoloContributo()/toString() methods ---
$ws.Range("J31").Value = "'+ calcoloContributo(): double"
$ws.Range("J32").Value = "'+ toString(): String"

# --- Professionale gains calcoloContributo()/toString() methods ---
$ws.Range("E34").Value = "'+ calcoloContributo(): double"
$ws.Range("E35").Value = "'+ toString(): String"

# --- Annotation notes with the actual contribution formulas, in the small
#     "Arial Unicode MS" 10pt note font used elsewhere for asides ---
$ws.Range("N8").Value = "125 * numeroStudenti + 9000 * numeroSediAggiuntive"
$ws.Range("N8").Font.Name = "Arial Unicode MS"
$ws.Range("N8").Font.Size = 10

$ws.Range("N16").Value = "150 * numeroStudenti + 1100 * numeroLaboratori + 9000 * numeroSediAggiuntive"
$ws.Range("N16").Font.Name = "Arial Unicode MS"
$ws.Range("N16").Font.Size = 10

$ws.Range("N30").Value = "3500 * numeroClassi + 6000 * numeroLaboratori"
$ws.Range("N30").Font.Name = "Arial Unicode MS"
$ws.Range("N30").Font.Size = 10

$ws.Range("A33").Value = "2400 * numeroClassi + 3000 * numeroLaboratori + contributoRegionale"
$ws.Range("A33").Font.Name = "Arial Unicode MS"
$ws.Range("A33").Font.Size = 10

# --- Move the active selection to A33, matching the author's last edit ---
$ws.Range("A33").Select() | Out-Null
